# Applies the "Updated symbol list" GitHub Actions commit to cryptos.xlsx
# Most edits are routine price/volume refreshes (column D "Price" and the
# rank-prefixed column E "Volume(1h)" label). Rows 12-20 additionally show a
# re-ranking: LiechtensteinCryptoassetsExchange jumped from row 20 up to row
# 12, shifting BitrueCoin..BitpandaEcosystemToken down by one row each, so
# those rows get full B/C/D/E rewrites instead of just a D/E refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple price (column D) refreshes, rows 2-11 ---
$ws.Range("D2").Value = '244.07'
$ws.Range("D4").Value = '5.186'
$ws.Range("D5").Value = '0.05736'
$ws.Range("D6").Value = '6.494'
$ws.Range("D7").Value = '3.116'
$ws.Range("D8").Value = '0.8092'
$ws.Range("D9").Value = '0.8392'
$ws.Range("D10").Value = '0.1339'
$ws.Range("D11").Value = '0.06959'

# --- Rows 12-20: re-ranking shift (full B/C/D/E rewrite per row) ---
# Row 12: -> LiechtensteinCryptoassetsExchange
$ws.Range("B12").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C12").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D12").Value = '0.03133'
$ws.Range("E12").Value = '11LiechtensteinCryptoassetsExchangeLCX'
# Row 13: -> BitrueCoin
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").Value = '0.02828'
$ws.Range("E13").Value = '12BitrueCoinBTR'
# Row 14: -> BitMartToken
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").Value = '0.09366'
$ws.Range("E14").Value = '13BitMartTokenBMX'
# Row 15: -> BitForexToken
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").Value = '0.001509'
$ws.Range("E15").Value = '14BitForexTokenBF'
# Row 16: -> One
$ws.Range("B16").Value = 'One'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D16").Value = '0.0005973'
$ws.Range("E16").Value = '15OneONE'
# Row 17: -> TigerCash
$ws.Range("B17").Value = 'TigerCash'
$ws.Range("C17").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D17").Value = '0.006233'
$ws.Range("E17").Value = '16TigerCashTCH'
# Row 18: -> LEO
$ws.Range("B18").Value = 'LEO'
$ws.Range("C18").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D18").Value = '3.503'
$ws.Range("E18").Value = '17LEOLEO'
# Row 19: -> BTSEToken
$ws.Range("B19").Value = 'BTSEToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D19").Value = '2.092'
$ws.Range("E19").Value = '18BTSETokenBTSE'
# Row 20: -> BitpandaEcosystemToken
$ws.Range("B20").Value = 'BitpandaEcosystemToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D20").Value = '0.3191'
$ws.Range("E20").Value = '19BitpandaEcosystemTokenBEST'

# --- Rows 22-28: price refreshes, plus UpBots label change on row 28 ---
$ws.Range("D22").Value = '3.744'
$ws.Range("D23").Value = '0.04655'
$ws.Range("D24").Value = '0.1330'
$ws.Range("D25").Value = '0.001235'
$ws.Range("D26").Value = '0.004264'
$ws.Range("D27").Value = '0.00008702'
$ws.Range("D28").Value = '0.0001987'
$ws.Range("E28").Value = '27UpBotsUBXT'

# --- Rows 40-50: price refreshes, plus Best/Worst-in-24h label churn ---
$ws.Range("D40").Value = '0.03611'
$ws.Range("D41").Value = '0.006311'
$ws.Range("E41").Value = '40KickTokenKICK'
$ws.Range("D42").Value = '0.1049'
$ws.Range("D43").Value = '0.002915'
$ws.Range("E43").Value = '42CEJICEJIBestin24h'
$ws.Range("D44").Value = '0.007337'
$ws.Range("E44").Value = '43LocalTradersLCTWorstin24h'
$ws.Range("D45").Value = '0.00005290'
$ws.Range("D47").Value = '0.3442'
$ws.Range("D48").Value = '0.002286'
$ws.Range("D49").Value = '0.00002101'
$ws.Range("D50").Value = '0.0002001'
